$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.616.76"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.823.75"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.61"
$ws.Range("E5").Value = "  +5.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.00"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  +5.82%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.42"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0852"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.01"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.272.83"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.827.46"
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.892"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.458.21"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("E19").Value = "  +7.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.16"
$ws.Range("E20").Value = "  -4.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.38"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0992"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.44"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.61"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("E25").Value = "  +2.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.65"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.13"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.65"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.81"
$ws.Range("E33").Value = "  +4.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0439"
$ws.Range("E34").Value = "  +24.36%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.95"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.06"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.05"
$ws.Range("E40").Value = "  -5.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.75"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "126.58"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.51"
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.30"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.078.92"
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.67"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.920"
$ws.Range("E50").Value = "  +5.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.85"
$ws.Range("E51").Value = "  +0.93%  "
